$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Fluent" -> "Experienced"  (Languages: ... section)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Fluent", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Experienced", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Two ListParagraph bullet items under the Gautrelet Scholarship /
#    AlgoExpert certification entries: spacing before=1/after=0 (twips)
#    -> before=57/after=57 (twips) i.e. 0.05pt -> 2.85pt
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Completed 100+ Software Engineering challenges", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$r1.Paragraphs.Item(1).SpaceBefore = 2.85
$r1.Paragraphs.Item(1).SpaceAfter = 2.85

$r2 = $d.Content
$r2.Find.Execute("Gautrelet", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$r2.Paragraphs.Item(1).SpaceBefore = 2.85
$r2.Paragraphs.Item(1).SpaceAfter = 2.85

# ------------------------------------------------------------------
# 3) "Al" <tab> "May" -> "Al" + spaces, tab removed (keeps two runs)
#    Target text: "...Mobile, Al[tab]May 2017-May 2018"
# ------------------------------------------------------------------
$tabFind = $d.Content
$found = $tabFind.Find.Execute("Al" + [char]9 + "May", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if ($found) {
    $segStart = $tabFind.Start
    $spaces = "".PadLeft(63)
    $alRange = $d.Range($segStart, $segStart + 2)
    $alRange.Text = "Al" + $spaces

    $tabRange = $d.Range($segStart + 2 + 63, $segStart + 3 + 63)
    $tabRange.Text = ""
}

# ------------------------------------------------------------------
# 4) "2017-May " -> "2017 - May " (spaces around en dash), scoped to
#    the area after the Al/May fix above to avoid the unrelated
#    "Spring Hill College" education entry earlier in the document.
# ------------------------------------------------------------------
$enDash = [char]8211
$searchText = "2017$($enDash)May "
$replaceText = "2017 $($enDash) May "
$tailStart = $d.Range(2000, 2001).Start
$tailRange = $d.Range($tailStart, $d.Content.End)
$tailRange.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2) | Out-Null

# ------------------------------------------------------------------
# 5) Add the built-in "Numbering Symbols" character style (materialized
#    by Word because of the Symbol-font bullet list used in the resume)
# ------------------------------------------------------------------
$newStyle = $d.Styles.Add("Numbering Symbols", 2)
$newStyle.QuickStyle = $true
